# The "Förändrad" (Changed) date column (C) for every data row is bumped
# by one day (46081 -> 46082), i.e. an automatic "last updated" stamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item(1, 1).Worksheet.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value -ne $null) {
        $cell.Value = 46082
    }
}
